# Reproduce the exact editing sequence used by the original author (this
# matters so the regenerated xl/sharedStrings.xml gets new <si> entries
# appended in the same order/index as the target diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: append the new "acquisition" item rows right after the
# existing last row (82), i.e. at rows 83-94.
$ws.Range("A83").Value = "acquisition_item_desktop"
$ws.Range("B83").Value = "Desktop"
$ws.Range("A84").Value = "acquisition_item_keyboard"
$ws.Range("B84").Value = "Keyboard"
$ws.Range("A85").Value = "acquisition_item_monitor"
$ws.Range("B85").Value = "Monitor"
$ws.Range("A86").Value = "acquisition_item_monitorCable"
$ws.Range("B86").Value = "Monitor Cable"
$ws.Range("A87").Value = "acquisition_item_mouse"
$ws.Range("B87").Value = "Mouse"
$ws.Range("A88").Value = "acquisition_item_networkCable"
$ws.Range("B88").Value = "Network Cable"
$ws.Range("A89").Value = "acquisition_item_powerCable"
$ws.Range("B89").Value = "Power Cable"
$ws.Range("A90").Value = "acquisition_item_usbStick"
$ws.Range("B90").Value = "USB Flash Drive"
$ws.Range("A91").Value = "acquisition_item_deskFan"
$ws.Range("B91").Value = "Desk Fan"
$ws.Range("A92").Value = "acquisition_item_photo"
$ws.Range("B92").Value = "Photo"
$ws.Range("A93").Value = "acquisition_item_stapler"
$ws.Range("B93").Value = "Stapler"
$ws.Range("A94").Value = "acquisition_title"
$ws.Range("B94").Value = "Acquisition"

# --- Step 2: insert a new row above row 13 for the "back" button key,
# pushing the chain_of_custody.. block (old rows 13-82 plus the rows just
# added above, now rows 13-94) down by one, to rows 14-95.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "back"
$ws.Range("B13").Value = "BACK"

# --- Step 3: append the final "Item(s):" label row at the very end.
$ws.Range("A96").Value = "acquisition_items"
$ws.Range("B96").Value = "Item(s):"

# --- Step 4: restore the view state (active cell) to match the freshly
# appended last row.
$ws.Range("A96").Select()
